$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 4, 5 and 6 each hold one species observation. The edit cyclically
# rotates the observation-specific columns "up" by one row and wraps the
# top one back to the bottom:
#   new row 4 <- old row 5
#   new row 5 <- old row 6
#   new row 6 <- old row 4
# (Id/Taxonsorteringsordning/species fields/count/unit/location/coords.)
$cols = @("A", "B", "D", "E", "F", "G", "H", "I", "J", "P", "Q", "R")

# Snapshot the "before" values for rows 4, 5 and 6 first (Value2 is the
# reliable read accessor in this host; plain .Value only works for writes).
$row4 = @{}
$row5 = @{}
$row6 = @{}
foreach ($col in $cols) {
    $row4[$col] = $ws.Range($col + "4").Value2
    $row5[$col] = $ws.Range($col + "5").Value2
    $row6[$col] = $ws.Range($col + "6").Value2
}

# Now write the rotated values back.
foreach ($col in $cols) {
    $ws.Range($col + "4").Value = $row5[$col]
    $ws.Range($col + "5").Value = $row6[$col]
    $ws.Range($col + "6").Value = $row4[$col]
}
